$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr1 = New-Object 'object[,]' 24,9
$arr1[0,0] = 14.06510505766926
$arr1[0,1] = 6.308894975215158
$arr1[0,2] = 3.272119936287464
$arr1[0,3] = 6.80733231627563
$arr1[0,4] = 34.12660784679228
$arr1[0,5] = 45.83776948347001
$arr1[0,6] = 4.257287850445218
$arr1[0,7] = 5.041259206834757
$arr1[0,8] = 14.61400242503018
$arr1[1,0] = 13.23134699367921
$arr1[1,1] = 6.06866967883084
$arr1[1,2] = 3.141402032509065
$arr1[1,3] = 6.613790033188543
$arr1[1,4] = 33.38249724865857
$arr1[1,5] = 44.68972946513802
$arr1[1,6] = 4.5032949467525
$arr1[1,7] = 5.246230834177719
$arr1[1,8] = 14.45672297786143
$arr1[2,0] = 12.69298975496817
$arr1[2,1] = 5.917433283690186
$arr1[2,2] = 3.060047583724399
$arr1[2,3] = 6.491194543686393
$arr1[2,4] = 32.91835749007371
$arr1[2,5] = 43.96975729220365
$arr1[2,6] = 4.659714024264876
$arr1[2,7] = 5.377009069297619
$arr1[2,8] = 14.35990581970669
$arr1[3,0] = 12.46487163668696
$arr1[3,1] = 5.859060751994606
$arr1[3,2] = 3.029895374995019
$arr1[3,3] = 6.43931496315681
$arr1[3,4] = 32.70758026583299
$arr1[3,5] = 43.63829123285038
$arr1[3,6] = 4.725328469859664
$arr1[3,7] = 5.433850616150563
$arr1[3,8] = 14.31362291048865
$arr1[4,0] = 12.4239000826298
$arr1[4,1] = 5.854485539176186
$arr1[4,2] = 3.028836217261077
$arr1[4,3] = 6.429419788136699
$arr1[4,4] = 32.64802928495252
$arr1[4,5] = 43.54074773740948
$arr1[4,6] = 4.736917282215979
$arr1[4,7] = 5.44626822882252
$arr1[4,8] = 14.29763882341915
$arr1[5,0] = 12.68269228163635
$arr1[5,1] = 5.930605865703779
$arr1[5,2] = 3.070291173606062
$arr1[5,3] = 6.487189237729495
$arr1[5,4] = 32.84926646288005
$arr1[5,5] = 43.85080118340656
$arr1[5,6] = 4.662237268752856
$arr1[5,7] = 5.385471149226329
$arr1[5,8] = 14.33673989041984
$arr1[6,0] = 13.77416134539979
$arr1[6,1] = 6.244662730907613
$arr1[6,2] = 3.240661586461362
$arr1[6,3] = 6.737190238910129
$arr1[6,4] = 33.7866025819976
$arr1[6,5] = 45.29871289456621
$arr1[6,6] = 4.343088120527143
$arr1[6,7] = 5.120658117489641
$arr1[6,8] = 14.53058941221476
$arr1[7,0] = 15.71284518057224
$arr1[7,1] = 6.80729577646801
$arr1[7,2] = 3.546883754407277
$arr1[7,3] = 7.202175849497329
$arr1[7,4] = 35.67310045304805
$arr1[7,5] = 48.20142497769632
$arr1[7,6] = 3.756048474494799
$arr1[7,7] = 4.626984728674308
$arr1[7,8] = 14.94719076538099
$arr1[8,0] = 16.98733976357979
$arr1[8,1] = 7.232904430857201
$arr1[8,2] = 3.762239287781317
$arr1[8,3] = 7.437607573241271
$arr1[8,4] = 36.68285826577396
$arr1[8,5] = 49.72435363195387
$arr1[8,6] = 3.376154666801476
$arr1[8,7] = 4.298487998824559
$arr1[8,8] = 15.14427530627029
$arr1[9,0] = 17.4144576583691
$arr1[9,1] = 7.678476017556322
$arr1[9,2] = 3.816283105549543
$arr1[9,3] = 6.855468934048872
$arr1[9,4] = 34.49415547706983
$arr1[9,5] = 46.27503389709238
$arr1[9,6] = 3.975979148853099
$arr1[9,7] = 4.229329756791367
$arr1[9,8] = 14.38360561120001
$arr1[10,0] = 17.52137132490764
$arr1[10,1] = 7.983635206315219
$arr1[10,2] = 3.800040860254384
$arr1[10,3] = 6.394847319921102
$arr1[10,4] = 32.5283587076438
$arr1[10,5] = 43.19860822757578
$arr1[10,6] = 5.048095394147619
$arr1[10,7] = 4.220485699198792
$arr1[10,8] = 13.74075956564114
$arr1[11,0] = 17.38324095160098
$arr1[11,1] = 8.219440716054047
$arr1[11,2] = 3.741888026688593
$arr1[11,3] = 6.006839606600063
$arr1[11,4] = 30.50519726871209
$arr1[11,5] = 40.03101245356423
$arr1[11,6] = 6.297600252310434
$arr1[11,7] = 4.269983197764249
$arr1[11,8] = 13.11627322891884
$arr1[12,0] = 17.17218120842997
$arr1[12,1] = 8.36325820059924
$arr1[12,2] = 3.683848853167134
$arr1[12,3] = 5.789864738671686
$arr1[12,4] = 29.05158085144935
$arr1[12,5] = 37.74857548594435
$arr1[12,6] = 7.22855379321889
$arr1[12,7] = 4.333074130844333
$arr1[12,8] = 12.68711403120014
$arr1[13,0] = 17.065593427717
$arr1[13,1] = 8.387014650229116
$arr1[13,2] = 3.664471624200087
$arr1[13,3] = 5.741705689190518
$arr1[13,4] = 28.66533593718645
$arr1[13,5] = 37.13869170633419
$arr1[13,6] = 7.452875497819007
$arr1[13,7] = 4.36510075785751
$arr1[13,8] = 12.5793760186124
$arr1[14,0] = 16.55930762879785
$arr1[14,1] = 8.182731832719096
$arr1[14,2] = 3.590008323051461
$arr1[14,3] = 5.70626824924527
$arr1[14,4] = 28.56137674595188
$arr1[14,5] = 36.97850774997434
$arr1[14,6] = 7.317527073234846
$arr1[14,7] = 4.500096124995821
$arr1[14,8] = 12.60504472108595
$arr1[15,0] = 16.28024192466454
$arr1[15,1] = 7.941716890523103
$arr1[15,2] = 3.564715737664871
$arr1[15,3] = 5.783373945873133
$arr1[15,4] = 29.25674538870685
$arr1[15,5] = 38.07575883039216
$arr1[15,6] = 6.700294160788587
$arr1[15,7] = 4.570963751689398
$arr1[15,8] = 12.85011191681394
$arr1[16,0] = 16.17593076587333
$arr1[16,1] = 7.641034321509391
$arr1[16,2] = 3.570117235452315
$arr1[16,3] = 6.029454899061316
$arr1[16,4] = 30.78734535188187
$arr1[16,5] = 40.49309948154198
$arr1[16,6] = 5.659894691007842
$arr1[16,7] = 4.583598250563234
$arr1[16,8] = 13.33854418337025
$arr1[17,0] = 16.20870215234283
$arr1[17,1] = 7.374500217726873
$arr1[17,2] = 3.6139082234033
$arr1[17,3] = 6.46774267227595
$arr1[17,4] = 32.80229744958862
$arr1[17,5] = 43.65493608792558
$arr1[17,6] = 4.52384773195081
$arr1[17,7] = 4.559741847163447
$arr1[17,8] = 13.96919116174252
$arr1[18,0] = 16.64275660086033
$arr1[18,1] = 7.162910541254057
$arr1[18,2] = 3.733656934027333
$arr1[18,3] = 7.364068025012534
$arr1[18,4] = 36.22299483849392
$arr1[18,5] = 48.99421883108078
$arr1[18,6] = 3.478997007153221
$arr1[18,7] = 4.408379422200412
$arr1[18,8] = 15.02329376664596
$arr1[19,0] = 17.60314488967096
$arr1[19,1] = 7.437324360766709
$arr1[19,2] = 3.898447314728777
$arr1[19,3] = 7.672108924356806
$arr1[19,4] = 37.48570578254825
$arr1[19,5] = 50.92507649809882
$arr1[19,6] = 3.152876673274757
$arr1[19,7] = 4.14633465158162
$arr1[19,8] = 15.33687889989804
$arr1[20,0] = 18.20635285180175
$arr1[20,1] = 7.611060985388321
$arr1[20,2] = 3.99055985012112
$arr1[20,3] = 7.829530121816125
$arr1[20,4] = 38.21986449531011
$arr1[20,5] = 52.05168315252219
$arr1[20,6] = 2.954911296435297
$arr1[20,7] = 3.972429419438719
$arr1[20,8] = 15.51829251817908
$arr1[21,0] = 17.89327491028746
$arr1[21,1] = 7.504457827707149
$arr1[21,2] = 3.931607339626766
$arr1[21,3] = 7.749136938614846
$arr1[21,4] = 37.89653665985571
$arr1[21,5] = 51.56717313831414
$arr1[21,6] = 3.058602313425717
$arr1[21,7] = 4.055205837516958
$arr1[21,8] = 15.44568490132298
$arr1[22,0] = 16.6478737750528
$arr1[22,1] = 7.11626462241958
$arr1[22,2] = 3.719818510481084
$arr1[22,3] = 7.431944599230935
$arr1[22,4] = 36.54621876713781
$arr1[22,5] = 49.51483236592602
$arr1[22,6] = 3.461493859038429
$arr1[22,7] = 4.39156527670162
$arr1[22,8] = 15.13208892106165
$arr1[23,0] = 15.20033225408193
$arr1[23,1] = 6.681569830642689
$arr1[23,2] = 3.482322428362618
$arr1[23,3] = 7.074382435571453
$arr1[23,4] = 35.05967636939739
$arr1[23,5] = 47.24148201868657
$arr1[23,6] = 3.9132375598701
$arr1[23,7] = 4.770191593901006
$arr1[23,8] = 14.79610357682358
$ws.Range("B2:J25").Value = $arr1

$arr2 = New-Object 'object[,]' 24,3
$arr2[0,0] = 5.956124619795937
$arr2[0,1] = 10.72932148134129
$arr2[0,2] = 5.999870472189118
$arr2[1,0] = 5.876855722801857
$arr2[1,1] = 10.05523597508855
$arr2[1,2] = 5.81611557747113
$arr2[2,0] = 5.826303039633578
$arr2[2,1] = 9.618460346679207
$arr2[2,2] = 5.701486074321997
$arr2[3,0] = 5.803976840250054
$arr2[3,1] = 9.43662685274152
$arr2[3,2] = 5.656544615903469
$arr2[4,0] = 5.798709994927463
$arr2[4,1] = 9.408419599797545
$arr2[4,2] = 5.651727738315972
$arr2[5,0] = 5.821857105027163
$arr2[5,1] = 9.622259863635071
$arr2[5,2] = 5.708090400212797
$arr2[6,0] = 5.923856088770586
$arr2[6,1] = 10.50928323230934
$arr2[6,2] = 5.946095383481208
$arr2[7,0] = 6.116058277264609
$arr2[7,1] = 12.05741759498994
$arr2[7,2] = 6.3844982216456
$arr2[8,0] = 6.206823480509083
$arr2[8,1] = 13.09645114597016
$arr2[8,2] = 6.634124096722892
$arr2[9,0] = 5.949497704563992
$arr2[9,1] = 13.63139060110016
$arr2[9,2] = 6.214423039412488
$arr2[10,0] = 5.78939785090925
$arr2[10,1] = 13.86737561176817
$arr2[10,2] = 5.808749471857796
$arr2[11,0] = 5.691890017430388
$arr2[11,1] = 13.90241021985043
$arr2[11,2] = 5.397163742277474
$arr2[12,0] = 5.663045319354099
$arr2[12,1] = 13.8384417877894
$arr2[12,2] = 5.111902329607644
$arr2[13,0] = 5.659339765295813
$arr2[13,1] = 13.78127476315311
$arr2[13,2] = 5.041770691612684
$arr2[14,0] = 5.634744015980669
$arr2[14,1] = 13.35734601099239
$arr2[14,2] = 5.011459008818604
$arr2[15,0] = 5.624105623595184
$arr2[15,1] = 13.06034025898363
$arr2[15,2] = 5.139591446124275
$arr2[16,0] = 5.665566945478751
$arr2[16,1] = 12.84130174938744
$arr2[16,2] = 5.42771423108984
$arr2[17,0] = 5.800511888621328
$arr2[17,1] = 12.71609290877447
$arr2[17,2] = 5.848996472484554
$arr2[18,0] = 6.170228135373891
$arr2[18,1] = 12.84660118376708
$arr2[18,2] = 6.585909066546527
$arr2[19,0] = 6.300891777600224
$arr2[19,1] = 13.59144919728455
$arr2[19,2] = 6.867373757252718
$arr2[20,0] = 6.369197410290926
$arr2[20,1] = 14.05766207282096
$arr2[20,2] = 7.008206909993604
$arr2[21,0] = 6.337202800652165
$arr2[21,1] = 13.80526627523457
$arr2[21,2] = 6.926139602924805
$arr2[22,0] = 6.20615345815286
$arr2[22,1] = 12.81811112719806
$arr2[22,2] = 6.620888154459569
$arr2[23,0] = 6.058668388815474
$arr2[23,1] = 11.66740324206201
$arr2[23,2] = 6.279758642624876
$ws.Range("L2:N25").Value = $arr2

